# Revert "added requirement openpyxl":
#  - drop the "sex" column from the test_file sheet (column E), shifting
#    customer_type left into its place
#  - append four more sample rows to Sheet2
#  - update sheet selections / active sheet to match the new session state

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("test_file")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Remove the "sex" column (E) from test_file; customer_type (old F) shifts to E.
$ws1.Columns.Item(5).Delete()

# Extend Sheet2 with four more rows following the existing 1,2,3 / 2,3,4 pattern.
$ws2.Range("A4").Value = 1
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 3
$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = 3
$ws2.Range("C5").Value = 4
$ws2.Range("A6").Value = 1
$ws2.Range("B6").Value = 2
$ws2.Range("C6").Value = 3
$ws2.Range("A7").Value = 2
$ws2.Range("B7").Value = 3
$ws2.Range("C7").Value = 4

# Update the view/selection state: Sheet2's selection moves to F16, and
# test_file becomes the active/selected tab with I8 selected.
$null = $ws2.Range("F16").Select()
$null = $ws1.Activate()
$null = $ws1.Range("I8").Select()
